$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Texas Notes" worksheet between "Calculations" and
#    "EoDSDwSP" (reverting the drop of the Texas-specific analysis tab).
# ---------------------------------------------------------------------------
$calcSheet = $wb.Worksheets.Item("Calculations")
$eodSheet  = $wb.Worksheets.Item("EoDSDwSP")
$txSheet   = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $calcSheet)
$txSheet.Name = "Texas Notes"

# Worksheet references captured before the insertion point shifted -
# re-resolve them by name now that the sheet collection has changed.
$calcSheet = $wb.Worksheets.Item("Calculations")
$eodSheet  = $wb.Worksheets.Item("EoDSDwSP")

# Column widths (characters)
$txSheet.Columns.Item(1).ColumnWidth = 22.166666666666668
$txSheet.Columns.Item(2).ColumnWidth = 67.5
$txSheet.Columns.Item(3).ColumnWidth = 67.5
$txSheet.Columns.Item(4).ColumnWidth = 19.166666666666668

# Data table header row (written first so the new shared strings line up
# with the order they were originally authored in)
$txSheet.Range("A5").Value = "Type"
$txSheet.Range("A5").Font.Bold = $true
$txSheet.Range("B5").Value = "BAU Deployment 2015-2022 (square pixels measured)"
$txSheet.Range("C5").Value = "Deployment with Extended ITC 2015-2022 (square pixels measured)"
$txSheet.Range("D5").Value = "% Increase due to ITC"
$txSheet.Range("B5:D5").Font.Bold = $true
$txSheet.Range("B5:D5").HorizontalAlignment = -4152

# Notes about the data source
$txSheet.Range("A1").Value = "The source has Texas specific data, but it's in graphical form. "
$txSheet.Range("A2").Value = "So, I used some visual editing software to measure the areas of the Texas specific data. This has some error associated with it."

# Data rows
$txSheet.Range("A6").Value = "Residential"
$txSheet.Range("B6").Value = 26477
$txSheet.Range("C6").Value = 29137
$txSheet.Range("D6").Formula = "=(C6-B6)/B6"
$txSheet.Range("D6").NumberFormat = "0.00%"

$txSheet.Range("A7").Value = "Commercial"
$txSheet.Range("B7").Value = 11010
$txSheet.Range("C7").Value = 20634
$txSheet.Range("D7").Formula = "=(C7-B7)/B7"
$txSheet.Range("D7").NumberFormat = "0.00%"

# ITC incentive level block
$txSheet.Range("A9").Value = "ITC Incentive Level"
$txSheet.Range("A9").Font.Bold = $true
$txSheet.Range("A9").HorizontalAlignment = -4152
$txSheet.Range("A10").Value = 0.3
$txSheet.Range("A10").NumberFormat = "0%"
$txSheet.Range("B10").Value = "of system cost"

# Elasticity block
$txSheet.Range("A12").Value = "Elasticity of Distributed Solar Deployment with respect to ITC Incentive Level"
$txSheet.Range("A12").Font.Bold = $true

$txSheet.Range("A13").Value = "Residential"
$txSheet.Range("B13").Formula = "=D6/`$A`$10"
$txSheet.Range("B13").NumberFormat = "0.000"

$txSheet.Range("A14").Value = "Commercial"
$txSheet.Range("B14").Formula = "=D7/`$A`$10"
$txSheet.Range("B14").NumberFormat = "0.000"

# Concluding notes
$txSheet.Range("A16").Value = "The main point here is that Texas residential is less elastic and commericial is more elastic than the national average. Because the measuring technique I used has some error in it, I will average the numbres"
$txSheet.Range("A17").Value = "above with the national numbers in the ""Calculations"" tab to come up with something a bit more conservative, in case my measuring error is high."

$txSheet.Range("A19").Value = "Residential"
$txSheet.Range("B19").Formula = "=AVERAGE(B13,Calculations!B9)"
$txSheet.Range("B19").NumberFormat = "0.000"

$txSheet.Range("A20").Value = "Commercial"
$txSheet.Range("B20").Formula = "=AVERAGE(B14,Calculations!B10)"
$txSheet.Range("B20").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# 2. Point EoDSDwSP's elasticity figures at the new Texas Notes averages
#    instead of straight at Calculations.
# ---------------------------------------------------------------------------
$eodSheet.Range("B2").Formula = "='Texas Notes'!B19"
$eodSheet.Range("B4").Formula = "='Texas Notes'!B20"

# ---------------------------------------------------------------------------
# 3. Add the source hyperlink on the About sheet.
# ---------------------------------------------------------------------------
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Hyperlinks.Add($aboutSheet.Range("B6"), "http://www.seia.org/sites/default/files/resources/BNEF_SEIA%20Solar%20Forecast_15%20September%202015.pdf") | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore the selections / active tab seen in the authored workbook.
# ---------------------------------------------------------------------------
$aboutSheet.Range("B6").Select() | Out-Null
$calcSheet.Range("B9").Select() | Out-Null
$txSheet.Range("B21").Select() | Out-Null
$eodSheet.Activate() | Out-Null
$eodSheet.Range("B5").Select() | Out-Null

Write-Host "done"
